# Apply updated TPM-based recalculation values to LR-pair table (NATMI output)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1242136666666667
$ws.Range("H2").Value = 0.372641
$ws.Range("I2").Value = 0.002558276231729912
$ws.Range("J2").Value = 0.002558276231729912
$ws.Range("M2").Value = 1021.934916333333
$ws.Range("N2").Value = 3065.804749
$ws.Range("O2").Value = 0.8026347959846111
$ws.Range("P2").Value = 0.802634795984611
$ws.Range("Q2").Value = 126.9382830524565
$ws.Range("R2").Value = 1142.444547472109
$ws.Range("S2").Value = 0.002053361521326818
$ws.Range("T2").Value = 0.002053361521326817

# Row 3
$ws.Range("G3").Value = 0.1242136666666667
$ws.Range("H3").Value = 0.372641
$ws.Range("I3").Value = 0.002558276231729912
$ws.Range("J3").Value = 0.002558276231729912
$ws.Range("O3").Value = 0.04931810976893385
$ws.Range("P3").Value = 0.04931810976893384
$ws.Range("Q3").Value = 7.799756762079223
$ws.Range("R3").Value = 70.19781085871301
$ws.Range("S3").Value = 0.0001261693480157103
$ws.Range("T3").Value = 0.0001261693480157102

# Row 4
$ws.Range("G4").Value = 0.1242136666666667
$ws.Range("H4").Value = 0.372641
$ws.Range("I4").Value = 0.002558276231729912
$ws.Range("J4").Value = 0.002558276231729912
$ws.Range("M4").Value = 187.139577
$ws.Range("N4").Value = 561.418731
$ws.Range("O4").Value = 0.1469807262726385
$ws.Range("P4").Value = 0.1469807262726385
$ws.Range("Q4").Value = 23.245293037619
$ws.Range("R4").Value = 209.207637338571
$ws.Range("S4").Value = 0.0003760172985456914
$ws.Range("T4").Value = 0.0003760172985456912

# Row 5
$ws.Range("G5").Value = 0.1242136666666667
$ws.Range("H5").Value = 0.372641
$ws.Range("I5").Value = 0.002558276231729912
$ws.Range("J5").Value = 0.002558276231729912
$ws.Range("M5").Value = 1.357726666666667
$ws.Range("N5").Value = 4.073180000000001
$ws.Range("O5").Value = 0.001066367973816652
$ws.Range("P5").Value = 0.001066367973816652
$ws.Range("Q5").Value = 0.1686482075977778
$ws.Range("R5").Value = 1.51783386838
$ws.Range("S5").Value = 0.000002728063841693126
$ws.Range("T5").Value = 0.000002728063841693125

# Row 6
$ws.Range("I6").Value = 0.4915413271548217
$ws.Range("J6").Value = 0.4915413271548217
$ws.Range("M6").Value = 1021.934916333333
$ws.Range("N6").Value = 3065.804749
$ws.Range("O6").Value = 0.8026347959846111
$ws.Range("P6").Value = 0.802634795984611
$ws.Range("Q6").Value = 24389.63054281554
$ws.Range("R6").Value = 219506.6748853399
$ws.Range("S6").Value = 0.3945281728389153
$ws.Range("T6").Value = 0.3945281728389152

# Row 7
$ws.Range("I7").Value = 0.4915413271548217
$ws.Range("J7").Value = 0.4915413271548217
$ws.Range("O7").Value = 0.04931810976893385
$ws.Range("P7").Value = 0.04931810976893384
$ws.Range("S7").Value = 0.02424188912858892
$ws.Range("T7").Value = 0.02424188912858892

# Row 8
$ws.Range("I8").Value = 0.4915413271548217
$ws.Range("J8").Value = 0.4915413271548217
$ws.Range("M8").Value = 187.139577
$ws.Range("N8").Value = 561.418731
$ws.Range("O8").Value = 0.1469807262726385
$ws.Range("P8").Value = 0.1469807262726385
$ws.Range("Q8").Value = 4466.297285687433
$ws.Range("R8").Value = 40196.6755711869
$ws.Range("S8").Value = 0.07224710125823232
$ws.Range("T8").Value = 0.07224710125823229

# Row 9
$ws.Range("I9").Value = 0.4915413271548217
$ws.Range("J9").Value = 0.4915413271548217
$ws.Range("M9").Value = 1.357726666666667
$ws.Range("N9").Value = 4.073180000000001
$ws.Range("O9").Value = 0.001066367973816652
$ws.Range("P9").Value = 0.001066367973816652
$ws.Range("Q9").Value = 32.40367977340667
$ws.Range("R9").Value = 291.6331179606601
$ws.Range("S9").Value = 0.0005241639290852353
$ws.Range("T9").Value = 0.0005241639290852351

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.3378266666666667
$ws.Range("H10").Value = 1.01348
$ws.Range("I10").Value = 0.006957800658901278
$ws.Range("J10").Value = 0.006957800658901277
$ws.Range("M10").Value = 1021.934916333333
$ws.Range("N10").Value = 3065.804749
$ws.Range("O10").Value = 0.8026347959846111
$ws.Range("P10").Value = 0.802634795984611
$ws.Range("Q10").Value = 345.2368663351689
$ws.Range("R10").Value = 3107.13179701652
$ws.Range("S10").Value = 0.00558457291235882
$ws.Range("T10").Value = 0.005584572912358818

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.3378266666666667
$ws.Range("H11").Value = 1.01348
$ws.Range("I11").Value = 0.006957800658901278
$ws.Range("J11").Value = 0.006957800658901277
$ws.Range("O11").Value = 0.04931810976893385
$ws.Range("P11").Value = 0.04931810976893384
$ws.Range("Q11").Value = 21.21317161351556
$ws.Range("R11").Value = 190.91854452164
$ws.Range("S11").Value = 0.0003431455766460535
$ws.Range("T11").Value = 0.0003431455766460534

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.3378266666666667
$ws.Range("H12").Value = 1.01348
$ws.Range("I12").Value = 0.006957800658901278
$ws.Range("J12").Value = 0.006957800658901277
$ws.Range("M12").Value = 187.139577
$ws.Range("N12").Value = 561.418731
$ws.Range("O12").Value = 0.1469807262726385
$ws.Range("P12").Value = 0.1469807262726385
$ws.Range("Q12").Value = 63.22073949932
$ws.Range("R12").Value = 568.98665549388
$ws.Range("S12").Value = 0.001022662594105553
$ws.Range("T12").Value = 0.001022662594105552

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.3378266666666667
$ws.Range("H13").Value = 1.01348
$ws.Range("I13").Value = 0.006957800658901278
$ws.Range("J13").Value = 0.006957800658901277
$ws.Range("M13").Value = 1.357726666666667
$ws.Range("N13").Value = 4.073180000000001
$ws.Range("O13").Value = 0.001066367973816652
$ws.Range("P13").Value = 0.001066367973816652
$ws.Range("Q13").Value = 0.4586762740444446
$ws.Range("R13").Value = 4.1280864664
$ws.Range("S13").Value = 0.000007419575790852722
$ws.Range("T13").Value = 0.000007419575790852719

# Row 14
$ws.Range("G14").Value = 24.22548766666667
$ws.Range("H14").Value = 72.676463
$ws.Range("I14").Value = 0.498942595954547
$ws.Range("J14").Value = 0.498942595954547
$ws.Range("M14").Value = 1021.934916333333
$ws.Range("N14").Value = 3065.804749
$ws.Range("O14").Value = 0.8026347959846111
$ws.Range("P14").Value = 0.802634795984611
$ws.Range("Q14").Value = 24756.8717117692
$ws.Range("R14").Value = 222811.8454059228
$ws.Range("S14").Value = 0.4004686887120101
$ws.Range("T14").Value = 0.40046868871201

# Row 15
$ws.Range("G15").Value = 24.22548766666667
$ws.Range("H15").Value = 72.676463
$ws.Range("I15").Value = 0.498942595954547
$ws.Range("J15").Value = 0.498942595954547
$ws.Range("O15").Value = 0.04931810976893385
$ws.Range("P15").Value = 0.04931810976893384
$ws.Range("Q15").Value = 1521.192605559373
$ws.Range("R15").Value = 13690.73345003436
$ws.Range("S15").Value = 0.02460690571568316
$ws.Range("T15").Value = 0.02460690571568315

# Row 16
$ws.Range("G16").Value = 24.22548766666667
$ws.Range("H16").Value = 72.676463
$ws.Range("I16").Value = 0.498942595954547
$ws.Range("J16").Value = 0.498942595954547
$ws.Range("M16").Value = 187.139577
$ws.Range("N16").Value = 561.418731
$ws.Range("O16").Value = 0.1469807262726385
$ws.Range("P16").Value = 0.1469807262726385
$ws.Range("Q16").Value = 4533.547514558717
$ws.Range("R16").Value = 40801.92763102845
$ws.Range("S16").Value = 0.07333494512175497
$ws.Range("T16").Value = 0.07333494512175494

# Row 17
$ws.Range("G17").Value = 24.22548766666667
$ws.Range("H17").Value = 72.676463
$ws.Range("I17").Value = 0.498942595954547
$ws.Range("J17").Value = 0.498942595954547
$ws.Range("M17").Value = 1.357726666666667
$ws.Range("N17").Value = 4.073180000000001
$ws.Range("O17").Value = 0.001066367973816652
$ws.Range("P17").Value = 0.001066367973816652
$ws.Range("Q17").Value = 32.89159061803778
$ws.Range("R17").Value = 296.0243155623401
$ws.Range("S17").Value = 0.0005320564050988708
$ws.Range("T17").Value = 0.0005320564050988705
